$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: rows 16-29, alternating JAIME / EDUARDO by ascending period (2302 -> 2308)
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico

$rows = @(
    @{ r=16; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2302"; f=55592; g=1051690 },
    @{ r=17; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2302"; f=53276; g=1148189 },
    @{ r=18; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2303"; f=55592; g=1051690 },
    @{ r=19; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2303"; f=53276; g=1148189 },
    @{ r=20; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2304"; f=55592; g=1051690 },
    @{ r=21; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2304"; f=53276; g=1148189 },
    @{ r=22; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2305"; f=55592; g=1051690 },
    @{ r=23; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2305"; f=53276; g=1148189 },
    @{ r=24; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2306"; f=55592; g=1051690 },
    @{ r=25; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2306"; f=53276; g=1148189 },
    @{ r=26; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2307"; f=55592; g=1051690 },
    @{ r=27; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2307"; f=53276; g=1148189 },
    @{ r=28; b="CC"; c="9239323";  d="JAIME RAFAEL BARRIOS GUTIERREZ"; e="2308"; f=42067; g=1051690 },
    @{ r=29; b="CC"; c="15038313"; d="EDUARDO JOSE VERGARA FLOREZ";    e="2308"; f=45928; g=1148189 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
}
